$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1887.0435
$ws.Range("I40").Value = 1212.5
$ws.Range("J40").Value = 2246.8
$ws.Range("K40").Value = 1212.5
$ws.Range("L40").Value = 2246.8
$ws.Range("M40").Value = -1037.5
$ws.Range("N40").Value = -2596.8

$ws.Range("H112").Value = 1745.7142
$ws.Range("J112").Value = 1824.1666
$ws.Range("L112").Value = 5472.4998
$ws.Range("N112").Value = -7688.4998

$ws.Range("H113").Value = 2575.1292
$ws.Range("I113").Value = 2665.25
$ws.Range("J113").Value = 2411.2727
$ws.Range("K113").Value = 2665.25
$ws.Range("L113").Value = 2411.2727
$ws.Range("M113").Value = 588.75
$ws.Range("N113").Value = -8919.2727

$ws.Range("H114").Value = 24761.715
$ws.Range("J114").Value = 24761.715
$ws.Range("L114").Value = 24761.715
$ws.Range("N114").Value = -33439.715

$ws.Range("H132").Value = 2279.5588
$ws.Range("I132").Value = 2478.75
$ws.Range("J132").Value = 1350
$ws.Range("K132").Value = 7436.25
$ws.Range("L132").Value = 4050
$ws.Range("M132").Value = -4906.25
$ws.Range("N132").Value = -9110

$ws.Range("H138").Value = 1715.06
$ws.Range("I138").Value = 754.1372699999999
$ws.Range("J138").Value = 2715.204
$ws.Range("K138").Value = 2262.41181
$ws.Range("L138").Value = 8145.612000000001
$ws.Range("M138").Value = 2877.58819
$ws.Range("N138").Value = -18425.612

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21402.143
$ws.Range("I32").Value = 20777.863
$ws.Range("K32").Value = 20777.863
$ws.Range("M32").Value = -20490.863

$ws.Range("H45").Value = 1932.2632
$ws.Range("I45").Value = 1213.4445
$ws.Range("K45").Value = 1213.4445
$ws.Range("M45").Value = -836.4445000000001

$ws.Range("H61").Value = 979.70966
$ws.Range("I61").Value = 819.24
$ws.Range("K61").Value = 819.24
$ws.Range("M61").Value = -607.24

$ws.Range("H122").Value = 1052.4286
$ws.Range("I122").Value = 728.4
$ws.Range("J122").Value = 1862.5
$ws.Range("K122").Value = 2185.2
$ws.Range("L122").Value = 5587.5
$ws.Range("M122").Value = 264.8000000000002
$ws.Range("N122").Value = -10487.5

$ws.Range("H132").Value = 1024.381
$ws.Range("I132").Value = 880.74286
$ws.Range("K132").Value = 2642.22858
$ws.Range("M132").Value = -112.22858

$ws.Range("H136").Value = 979.70966
$ws.Range("I136").Value = 819.24
$ws.Range("K136").Value = 2457.72
$ws.Range("M136").Value = 92.27999999999975

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("I134").Value = 1529.7959
$ws.Range("J134").Value = 114195.445
$ws.Range("K134").Value = 4589.3877
$ws.Range("L134").Value = 342586.335
$ws.Range("M134").Value = -2054.3877
$ws.Range("N134").Value = -347656.335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2879.2
$ws.Range("I31").Value = 2648.9333
$ws.Range("J31").Value = 3570
$ws.Range("K31").Value = 2648.9333
$ws.Range("L31").Value = 3570
$ws.Range("M31").Value = -2353.9333
$ws.Range("N31").Value = -4160

$ws.Range("H34").Value = 2879.2
$ws.Range("I34").Value = 2648.9333
$ws.Range("J34").Value = 3570
$ws.Range("K34").Value = 2648.9333
$ws.Range("L34").Value = 3570
$ws.Range("M34").Value = -2446.9333
$ws.Range("N34").Value = -3974

$ws.Range("H99").Value = 2626.1538
$ws.Range("I99").Value = 2227.6191
$ws.Range("J99").Value = 4300
$ws.Range("K99").Value = 2227.6191
$ws.Range("L99").Value = 4300
$ws.Range("M99").Value = -729.6190999999999
$ws.Range("N99").Value = -7296

$ws.Range("H100").Value = 35780
$ws.Range("J100").Value = 35780
$ws.Range("L100").Value = 35780
$ws.Range("N100").Value = -37944

$ws.Range("H126").Value = 2626.1538
$ws.Range("I126").Value = 2227.6191
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 6682.8573
$ws.Range("L126").Value = 12900
$ws.Range("M126").Value = -4212.8573
$ws.Range("N126").Value = -17840

$ws.Range("H134").Value = 1990.24
$ws.Range("I134").Value = 2037.7646
$ws.Range("J134").Value = 1889.25
$ws.Range("K134").Value = 6113.293799999999
$ws.Range("L134").Value = 5667.75
$ws.Range("M134").Value = -3578.293799999999
$ws.Range("N134").Value = -10737.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H122").Value = 556.0222
$ws.Range("J122").Value = 689.0645
$ws.Range("L122").Value = 6201.5805
$ws.Range("N122").Value = -11101.5805

$ws.Range("H134").Value = 3018.2068
$ws.Range("I134").Value = 1853.3914
$ws.Range("K134").Value = 5560.174199999999
$ws.Range("M134").Value = -490.1741999999995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 24000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 24000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 24000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -25018

$ws.Range("H70").Value = 4079.5454
$ws.Range("J70").Value = 4233.3335
$ws.Range("L70").Value = 4233.3335
$ws.Range("N70").Value = -4773.3335

$ws.Range("H73").Value = 4079.5454
$ws.Range("J73").Value = 4233.3335
$ws.Range("L73").Value = 4233.3335
$ws.Range("N73").Value = -6105.3335

$ws.Range("H102").Value = 1432
$ws.Range("I102").Value = 1434.909
$ws.Range("J102").Value = 1400
$ws.Range("K102").Value = 1434.909
$ws.Range("L102").Value = 1400
$ws.Range("M102").Value = 187.0909999999999
$ws.Range("N102").Value = -4644

$ws.Range("H126").Value = 2797.2354
$ws.Range("I126").Value = 2752.7856
$ws.Range("J126").Value = 3004.6667
$ws.Range("K126").Value = 8258.356800000001
$ws.Range("L126").Value = 9014.000100000001
$ws.Range("M126").Value = -5788.356800000001
$ws.Range("N126").Value = -13954.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6176150
$ws.Range("I7").Value = 3598.7693
$ws.Range("J7").Value = 22224782
$ws.Range("K7").Value = 3598.7693
$ws.Range("L7").Value = 22224782
$ws.Range("M7").Value = -3486.7693
$ws.Range("N7").Value = -22225006

$ws.Range("H40").Value = 1123901.6
$ws.Range("I40").Value = 1684768.4
$ws.Range("J40").Value = 2168.3333
$ws.Range("K40").Value = 1684768.4
$ws.Range("L40").Value = 2168.3333
$ws.Range("M40").Value = -1684632.4
$ws.Range("N40").Value = -2440.3333

$ws.Range("H126").Value = 6176150
$ws.Range("I126").Value = 3598.7693
$ws.Range("J126").Value = 22224782
$ws.Range("K126").Value = 10796.3079
$ws.Range("L126").Value = 66674346
$ws.Range("M126").Value = -8326.3079
$ws.Range("N126").Value = -66679286

$ws.Range("H132").Value = 1817.6229
$ws.Range("I132").Value = 1655.8914
$ws.Range("J132").Value = 2313.6
$ws.Range("K132").Value = 4967.674199999999
$ws.Range("L132").Value = 6940.799999999999
$ws.Range("M132").Value = -2437.674199999999
$ws.Range("N132").Value = -12000.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1027.7142
$ws.Range("I126").Value = 1097
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 3291
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -821
$ws.Range("N126").Value = -7940

$ws.Range("H132").Value = 672.4143
$ws.Range("I132").Value = 568.62744
$ws.Range("J132").Value = 951
$ws.Range("K132").Value = 1705.88232
$ws.Range("L132").Value = 2853
$ws.Range("M132").Value = 824.1176800000001
$ws.Range("N132").Value = -7913
